$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the target range as Text so the numeric-looking strings
# (prices, percentages) are stored as literal text, matching the source data.
$ws.Range("D2:E50").NumberFormat = "@"

$ws.Range("D2").Value = "292.47"
$ws.Range("E2").Value = "-1.46%"
$ws.Range("D3").Value = "40.47"
$ws.Range("E3").Value = "-0.52%"
$ws.Range("D4").Value = "5.013"
$ws.Range("E4").Value = "-0.58%"
$ws.Range("D5").Value = "0.07298"
$ws.Range("E5").Value = "-1.67%"
$ws.Range("D6").Value = "4.279"
$ws.Range("E6").Value = "-0.76%"
$ws.Range("D7").Value = "1.557"
$ws.Range("E7").Value = "-1.61%"
$ws.Range("D8").Value = "0.9285"
$ws.Range("E8").Value = "0.38%"
$ws.Range("D10").Value = "0.1162"
$ws.Range("E10").Value = "-0.21%"
$ws.Range("D11").Value = "0.1752"
$ws.Range("E11").Value = "0.25%"
$ws.Range("D12").Value = "0.04368"
$ws.Range("E12").Value = "4.21%"
$ws.Range("D13").Value = "0.08679"
$ws.Range("E13").Value = "-1.00%"
$ws.Range("D14").Value = "0.1053"
$ws.Range("E14").Value = "0.04%"
$ws.Range("D15").Value = "0.001264"
$ws.Range("E15").Value = "-0.10%"
$ws.Range("D16").Value = "0.006021"
$ws.Range("E16").Value = "0.26%"
$ws.Range("D17").Value = "3.344"
$ws.Range("E17").Value = "-0.53%"
$ws.Range("D19").Value = "7.904"
$ws.Range("E19").Value = "3.42%"
$ws.Range("D20").Value = "0.1390"
$ws.Range("E20").Value = "2.17%"
$ws.Range("D21").Value = "0.2772"
$ws.Range("E21").Value = "-1.78%"
$ws.Range("D22").Value = "0.03921"
$ws.Range("E22").Value = "1.40%"
$ws.Range("E23").Value = "-2.16%"
$ws.Range("D24").Value = "0.003683"
$ws.Range("E24").Value = "1.42%"
$ws.Range("E25").Value = "-8.13%"
$ws.Range("D26").Value = "0.0003726"
$ws.Range("E26").Value = "-0.48%"
$ws.Range("D38").Value = "0.02314"
$ws.Range("E38").Value = "-0.03%"
$ws.Range("D39").Value = "0.05071"
$ws.Range("E39").Value = "1.21%"
$ws.Range("D40").Value = "0.005716"
$ws.Range("E40").Value = "37.11%"
$ws.Range("D41").Value = "0.007853"
$ws.Range("E41").Value = "1.45%"
$ws.Range("D42").Value = "0.1286"
$ws.Range("E42").Value = "0.67%"
$ws.Range("D43").Value = "0.007373"
$ws.Range("E43").Value = "-0.80%"
$ws.Range("D44").Value = "0.007245"
$ws.Range("E44").Value = "-7.45%"
$ws.Range("D45").Value = "0.2911"
$ws.Range("E45").Value = "-9.64%"
$ws.Range("D46").Value = "0.00006204"
$ws.Range("E46").Value = "-4.41%"
$ws.Range("E47").Value = "-0.47%"
$ws.Range("D48").Value = "0.04875"
$ws.Range("E48").Value = "-80.64%"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").Value = "-0.47%"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").Value = "-0.47%"

# Remove the temporary Text number format so the cells retain their original
# (default) style, matching the workbook before this edit.
$ws.Range("D2:E50").ClearFormats()
